$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00351604194135744
$ws.Range("C2").Value = 0.00342186224649965
$ws.Range("D2").Value = 0.00257424499277956
$ws.Range("E2").Value = 0.00838199284234319
$ws.Range("F2").Value = 0.965655804608526
$ws.Range("G2").Value = 0.990833176367175
$ws.Range("H2").Value = 0.00021975262133484
$ws.Range("I2").Value = 0.00351604194135744
$ws.Range("J2").Value = 0.00100458341181641
$ws.Range("K2").Value = 0.989985559113455
$ws.Range("L2").Value = 0.0896904627362341
$ws.Range("M2").Value = 0.996421171595404
$ws.Range("N2").Value = 0.000816224022100835
$ws.Range("O2").Value = 0.0181766811075532
$ws.Range("P2").Value = 0.000251145852954103
$ws.Range("Q2").Value = 0.997300182080743
$ws.Range("R2").Value = 0.000470898474288943
$ws.Range("S2").Value = 0.00508570352232059
$ws.Range("T2").Value = 0.00929239655930181
$ws.Range("U2").Value = 0.976706222138507
$ws.Range("V2").Value = 0.00831920637910467
$ws.Range("W2").Value = 0.00141269542286683
$ws.Range("X2").Value = 0.140861430275633
$ws.Range("B3").Value = 0.000502291705908206
$ws.Range("C3").Value = 0.0000627864632385258
$ws.Range("D3").Value = 0.000753437558862309
$ws.Range("E3").Value = 0.989608840334024
$ws.Range("F3").Value = 0.0113643498461732
$ws.Range("G3").Value = 0.00021975262133484
$ws.Range("H3").Value = 0.00021975262133484
$ws.Range("I3").Value = 0.0460538707854587
$ws.Range("J3").Value = 0.989483267407547
$ws.Range("K3").Value = 0.00021975262133484
$ws.Range("L3").Value = 0.000376718779431155
$ws.Range("M3").Value = 0.0000313932316192629
$ws.Range("N3").Value = 0.000408112011050418
$ws.Range("O3").Value = 0.0414704589690463
$ws.Range("P3").Value = 0.000784830790481572
$ws.Range("Q3").Value = 0.000282539084573366
$ws.Range("R3").Value = 0.458403968104477
$ws.Range("S3").Value = 0.52131600426948
$ws.Range("T3").Value = 0.98295347523074
$ws.Range("U3").Value = 0.00976329503359076
$ws.Range("V3").Value = 0.00100458341181641
$ws.Range("W3").Value = 0.010862058140265
$ws.Range("X3").Value = 0.00897846424310919
$ws.Range("B4").Value = 0.987788032900107
$ws.Range("C4").Value = 0.988007785521442
$ws.Range("D4").Value = 0.993344634896716
$ws.Range("E4").Value = 0.00113015633829346
$ws.Range("F4").Value = 0.0129967978903748
$ws.Range("G4").Value = 0.00734601619890752
$ws.Range("H4").Value = 0.991304074841464
$ws.Range("I4").Value = 0.00813084698938909
$ws.Range("J4").Value = 0.00885289131663213
$ws.Range("K4").Value = 0.00891567777987066
$ws.Range("L4").Value = 0.888930746531048
$ws.Range("M4").Value = 0.00128712249638978
$ws.Range("N4").Value = 0.00276260438249513
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.0000627864632385258
$ws.Range("Q4").Value = 0.000282539084573366
$ws.Range("R4").Value = 0.00043950524266968
$ws.Range("S4").Value = 0.42292961637471
$ws.Range("T4").Value = 0.00282539084573366
$ws.Range("U4").Value = 0.00361022163621523
$ws.Range("V4").Value = 0.986469517172098
$ws.Range("W4").Value = 0.987474100583914
$ws.Range("X4").Value = 0.816286808564074
$ws.Range("B5").Value = 0.00775412820995793
$ws.Range("C5").Value = 0.00800527406291204
$ws.Range("D5").Value = 0.00263703145601808
$ws.Range("E5").Value = 0.000156966158096314
$ws.Range("F5").Value = 0.00191498712877504
$ws.Range("G5").Value = 0.000941796948577887
$ws.Range("H5").Value = 0.00813084698938909
$ws.Range("I5").Value = 0.931908080617819
$ws.Range("J5").Value = 0.000188359389715577
$ws.Range("K5").Value = 0.000565078169146732
$ws.Range("L5").Value = 0.00339046901488039
$ws.Range("M5").Value = 0.000376718779431155
$ws.Range("N5").Value = 0.995479374646826
$ws.Range("O5").Value = 0.930526778426571
$ws.Range("P5").Value = 0.998461731650656
$ws.Range("Q5").Value = 0.00153826834934388
$ws.Range("R5").Value = 0.513938594838953
$ws.Range("S5").Value = 0.0367300809945376
$ws.Range("T5").Value = 0.00326489608840334
$ws.Range("U5").Value = 0.00662397187166447
$ws.Range("V5").Value = 0.00320210962516481
$ws.Range("W5").Value = 0.0000941796948577887
$ws.Range("X5").Value = 0.00191498712877504
